$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 24 data (new test case: "Analyse Module de Young") ---
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "Analyse Module de Young"
$ws.Range("D24").Value = "./Data/Exemple_WB100kN_Young"
$ws.Range("E24").Value = "Configurer l'échantillon comme cela : `nMesure extenso (W+B 100kN)`nMode d'analyse: Flexion 4pts`nGéométrie: Ronde`nD0 = 2`nL0 = 50`nL1= 20`nF_Max = 300`nF_Min = 150"
$ws.Range("F24").Value = "Résultats Correctes (Interface Shimadzu): A OBTENIR`nF_Max =  N`nAllong = mm`nRe = MPa`nRm =  Mpa`nDefo =  % `nE =  GPa"

# Match style of neighbouring rows: C/D/E wrap text, B centered/rotated, F red+wrap
$ws.Range("B24").Style = $ws.Range("B23").Style
$ws.Range("C24").Style = $ws.Range("C23").Style
$ws.Range("D24").Style = $ws.Range("D23").Style
$ws.Range("E24").Style = $ws.Range("E23").Style
$ws.Range("F24").Style = $ws.Range("F23").Style

# Row height for the new row
$ws.Rows.Item(24).RowHeight = 129.6

# Row 21 explicitly keeps its existing height, now marked as a custom height
$ws.Rows.Item(21).RowHeight = 115.2

# --- Extend the merged cell covering the "Comportement de l'analyse" block ---
$ws.Range("B21:B23").UnMerge()
$ws.Range("B21:B24").Merge()

# --- Update sheet view: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("A22").Select()
$ws.Range("F28").Select()
